$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row: "..._old" -> "..._FV2310" and "..._new" -> "..._FV2404"
#    (column K, "diff", is left untouched)
# ---------------------------------------------------------------------------
$oldHeaders = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")
$newHeaders = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $oldHeaders[$i]
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into a table so it gets an AutoFilter + Table1 definition
# ---------------------------------------------------------------------------
$range = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, "Table1", [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, keep focus on the bottom pane)
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
